# Annual Vehicle Maint Cost.xlsx - "Rail and aviation edits"
#
# Adds a new "rail cars per locomotive" input row (row 58) to the
# 'Cost Data' sheet (pushing the freight-rail block down by one row),
# and updates the per-car-per-year freight-rail cost calculation on the
# 'AVMC-passenger' summary sheet to divide the total freight-rail cost
# by the new "rail cars per locomotive" figure.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. 'Cost Data' sheet: insert a new row 58 with the rail-cars input.
# ---------------------------------------------------------------------
$costData = $wb.Worksheets.Item("Cost Data")

# Inserting the row shifts every cell/row below it down by one, and
# Excel automatically rewrites every formula (on this sheet and others)
# that referenced the shifted rows.
$costData.Rows("58:58").Insert()

$costData.Range("A58").Value = "rail cars per locomotive"
$costData.Range("B58").Value = 10
$costData.Range("B58").NumberFormat = "0.0"

# The hyperlink that used to sit on A66 is now one row down, on A67.
# (Row-insert doesn't renumber hyperlink anchors automatically, so the
# three external hyperlinks on this sheet are rebuilt explicitly.)
$costData.Hyperlinks.Delete()
$costData.Hyperlinks.Add($costData.Range("B7"), "https://www.energy.gov/eere/electricvehicles/electric-car-safety-maintenance-and-battery-life") | Out-Null
$costData.Hyperlinks.Add($costData.Range("B13"), "https://www.energy.gov/sites/prod/files/2014/03/f10/fuel_cell_mhe_cost.pdf") | Out-Null
$costData.Hyperlinks.Add($costData.Range("A67"), "https://www.insurance.com/motorcycle/is-riding-a-motorcycle-cheaper.html") | Out-Null

$costData.Range("B58").Select()

# ---------------------------------------------------------------------
# 2. 'AVMC-passenger' sheet: rail row (row 5) now averages the total
#    freight-rail maintenance cost per rail car, by dividing by the new
#    "rail cars per locomotive" figure on 'Cost Data'!$B$58.
# ---------------------------------------------------------------------
$passenger = $wb.Worksheets.Item("AVMC-passenger")

$passenger.Range("B5").Formula = "='Cost Data'!`$C90/'Cost Data'!`$B`$58"
$passenger.Range("C5").Formula = "='Cost Data'!`$B90/'Cost Data'!`$B`$58"
$passenger.Range("D5").Formula = "='Cost Data'!`$B90/'Cost Data'!`$B`$58"
$passenger.Range("E5").Formula = "='Cost Data'!`$B90/'Cost Data'!`$B`$58"
$passenger.Range("F5").Formula = "='Cost Data'!`$B90/'Cost Data'!`$B`$58"
$passenger.Range("G5").Formula = "='Cost Data'!`$B90/'Cost Data'!`$B`$58"
$passenger.Range("H5").Formula = "='Cost Data'!`$C90/'Cost Data'!`$B`$58"

$passenger.Range("H6").Select()

# ---------------------------------------------------------------------
# Restore the originally active sheet/tab so the workbook reopens the
# same way it did before these edits.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("About").Activate()
